$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing G / M / N values (rows 2-12) ---
$ws.Range("G2").Value = 3511.8525
$ws.Range("M2").Value = 43831.949999999997
$ws.Range("N2").Value = 55403.525000000001

$ws.Range("G3").Value = 34492.970000000001
$ws.Range("M3").Value = 609.47500000000002
$ws.Range("N3").Value = 706.96249999999998

$ws.Range("G4").Value = 18967.341250000001
$ws.Range("M4").Value = 275.3125
$ws.Range("N4").Value = 312.53750000000002

$ws.Range("G5").Value = 63839.386250000003
$ws.Range("M5").Value = 16886.900000000001
$ws.Range("N5").Value = 18954.924999999999

$ws.Range("G6").Value = 37088.654999999999
$ws.Range("M6").Value = 19110.849999999999
$ws.Range("N6").Value = 21103.337500000001

$ws.Range("G7").Value = 32566.193749999999
$ws.Range("M7").Value = 15017.6
$ws.Range("N7").Value = 15408.325000000001

$ws.Range("G8").Value = 17404.142500000002
$ws.Range("M8").Value = 7497.2749999999996
$ws.Range("N8").Value = 8292.1499999999996

$ws.Range("G9").Value = 15673.893749999999
$ws.Range("M9").Value = 2407.625
$ws.Range("N9").Value = 2611.0999999999999

$ws.Range("G10").Value = 20031.605
$ws.Range("M10").Value = 2763.5250000000001
$ws.Range("N10").Value = 3248.5875000000001

$ws.Range("G11").Value = 5611.2600000000002
$ws.Range("M11").Value = 15740.387500000001
$ws.Range("N11").Value = 17926.424999999999

$ws.Range("G12").Value = 13268.577499999999
$ws.Range("M12").Value = 13273.575000000001
$ws.Range("N12").Value = 14464.9

# --- Add new row 13: duplicate of row 4 with an updated path in column A ---
$ws.Range("A4:N4").Copy()
$ws.Range("A13").PasteSpecial(-4104)
$ws.Range("A13").Value = "C:\Users\rektplorer64\OneDrive\Pictures\Wallpapers\IMG_3037-scaled.jpg"

# --- Column A width now needs to fit the new, longer path string ---
$ws.Columns("A").ColumnWidth = 66.66666666666667
